# Insert a new data row at row 469 ("Vega Monumental Concepción" - Palta sheet).
# This shifts the previous rows 469-572 down to 470-573 (dimension grows from
# A1:T572 to A1:T573), and the newly inserted row 469 is populated with a new
# price record (Negra de La Cruz / Primera, from Región de O'Higgins).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 469 (and everything below it) down by one row.
$ws.Rows.Item(469).Insert()

# Populate the newly inserted, now-blank row 469 with the new record.
$ws.Range("A469").Value2 = 11
$ws.Range("B469").Value2 = "Vega Monumental Concepción"
$ws.Range("C469").Value2 = "Bíobío"
$ws.Range("D469").Value2 = 44722
$ws.Range("D469").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E469").Value2 = 8
$ws.Range("F469").Value2 = "Fruta"
$ws.Range("G469").Value2 = 100106
$ws.Range("H469").Value2 = "Oleaginosos"
$ws.Range("I469").Value2 = 100106002
$ws.Range("J469").Value2 = "Palta"
$ws.Range("K469").Value2 = "Negra de La Cruz"
$ws.Range("L469").Value2 = "Primera"
$ws.Range("M469").Value2 = 350
$ws.Range("N469").Value2 = 1300
$ws.Range("O469").Value2 = 1500
$ws.Range("P469").Value2 = 1414
$ws.Range("Q469").Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R469").Value2 = "Región de O'Higgins"
$ws.Range("S469").Value2 = 1414
$ws.Range("T469").Value2 = 1
